# Regenerate orders with updated distance/size codes.
#
# The experiment's distance and size condition labels were renumbered:
#   D51 -> D55
#   D64 -> D69
#   D80 -> D86
#   S30 -> S31
#
# These codes appear throughout the sheet (Condition, Filename_Left,
# Filename_Right, Distance, Size columns), always as a whole token inside
# underscore-delimited strings (e.g. "Face10_D51_S20", "Fixation_D51_l.png").
# A straightforward global text substitution across the used range
# reproduces the regenerated order file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

$rng.Replace("D51", "D55") | Out-Null
$rng.Replace("D64", "D69") | Out-Null
$rng.Replace("D80", "D86") | Out-Null
$rng.Replace("S30", "S31") | Out-Null
